$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy G1's formatting (font/border/alignment) onto H1 so the new header
# cell matches the other header cells (bold, bordered, centered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Header for new "Save" column
$ws.Range("H1").Value = "Save"

# New "Save" column values (1 = saved, 0 = not saved)
$values = @(1, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
